$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new row 10 (paper #5) -------------------------------------------
# Copy formatting from row 9 first (matches the look of the table), then
# fix up B10 (uses the plain style instead of the colored title style) and
# drop G10 (no value/format in this row).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C9").Copy()
$ws.Range("B10:F10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H9").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(10).RowHeight = 155.25

# Values - set in the same order the source workbook grew its shared
# string table so new si indices line up with the reference edit.
$ws.Range("B10").Value = "Intelligent Resource Scheduling for 5G Radio" + [char]10 + "Access Network Slicing"
$ws.Range("C10").Value = "propose an intelligent resource scheduling strategy (iRSS) for 5G RAN slicing. The main idea of iRSS is to exploit a collaborative learning framework which consists of deep learning (DL) in conjunction with Reinforcement Learning" + [char]10 + "(RL)"
$ws.Range("H10").Value = "DL is used to perform large time-scale resource" + [char]10 + "allocation, while RL is used to perform on-line resource scheduling for tackling small time-scale network dynamics, including" + [char]10 + "inaccurate prediction and unexpected network states"
$ws.Range("D10").Value = "There is two problem. 1: resource allocation of large time scale 2. RA of small scale, In the first step, we want to minimize the mean-square-error (MSE) between the predicted value of assigned resource block  and the actual traffic volume of it. In the second step we want to minimize assigned RB subject to the thresholds we need"
$ws.Range("E10").Value = "predicted value of assigned RB (pre assignement)- large time scale, fro small time scale-  number of slices and the time interval"
$ws.Range("F10").Value = "LTSM -large time scale, Q-learning, classic AC and HRSA algorithms for small time scale"
$ws.Range("A10").Value = 5

# --- Update row 9 (paper #4) text -----------------------------------------
$ws.Range("F9").Value = "first the modified deep deterministic" + [char]10 + "policy gradient (DDPG) for lower level and then the double deep-Q-network algorithm for upper level "
$ws.Range("E9").Value = "power, RB (find best policy for upper and lower level to assigne resources to slices ) for lower level and the guarantee bit rate and maximum rate for upper level"

# --- Sheet view: scroll position, zoom and selection ----------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F6").Select()
$excel.ActiveWindow.Zoom = 100
